$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: Condense the three detailed CORE COMPETENCIES paragraphs into
# a single summary paragraph: "Product Management & Strategy • Technical
# Product Development • Platform & Infrastructure"
# -------------------------------------------------------------------------

$bullet = [char]0x2022

$coreCompHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "CORE COMPETENCIES") {
        $coreCompHeading = $i
        break
    }
}

$firstDetail = $coreCompHeading + 1
$secondDetail = $coreCompHeading + 2
$thirdDetail = $coreCompHeading + 3

# Replace the first detail paragraph's text with the condensed summary,
# keeping its own paragraph mark (so no extra paragraph is introduced).
$p1 = $d.Paragraphs.Item($firstDetail)
$p1.Range.Text = "Product Management & Strategy $bullet Technical Product Development $bullet Platform & Infrastructure"

# Remove the now-redundant second and third detail paragraphs entirely.
$p2 = $d.Paragraphs.Item($secondDetail)
$p3 = $d.Paragraphs.Item($thirdDetail)
$removeRange = $d.Range($p2.Range.Start, $p3.Range.End)
$removeRange.Delete()

# -------------------------------------------------------------------------
# Change 2: Insert a new "TECHNICAL SKILLS" section (one Heading2 title and
# three body paragraphs) right before the closing "For a more detailed..."
# paragraph.
# -------------------------------------------------------------------------

$closingIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("For a more detailed")) {
        $closingIndex = $i
        break
    }
}

$precedingPara = $d.Paragraphs.Item($closingIndex - 1)
$insertionPoint = $precedingPara.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($closingIndex)
$newRange = $newPara.Range.Duplicate
$newRange.Collapse(1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>TECHNICAL SKILLS</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>PRODUCT MANAGEMENT &amp; STRATEGY Product Conception &amp; Ideation; Product Architecture &amp; Design; Product Lifecycle Management; B2B SaaS Development</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; API Development</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>PLATFORM &amp; INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Security &amp; Compliance</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xml) | Out-Null
